$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 429, shifting existing rows (429+) down by one.
$ws.Rows.Item(429).Insert()

# Populate the new row with the "suggestion" transaction status strings.
# (Order matches the shared-strings table insertion order: pt, key, en.)
$ws.Range("B429").Value = "Sugestão"
$ws.Range("A429").Value = "transaction_suggestion"
$ws.Range("C429").Value = "Suggestion"

# Match the row height used by the other simple (single-line) rows in the
# table; the headless engine does not compute real text-layout autofit, so
# this is set explicitly to 17 (same as the neighbouring rows).
$ws.Rows.Item(429).RowHeight = 17

# The "i18n" table needs to grow by one row to cover the newly inserted row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G567"))

# Update the view (top-left cell / active selection) as recorded after the edit.
$ws.Application.ActiveWindow.ScrollRow = 421
$ws.Range("C429").Select()
